$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..31 (columns A..R), per the diff:
# - rows 2-29 get their values reshuffled/updated to match the new weekly snapshot
# - rows 30-31 are brand new rows appended at the bottom
$data = @(
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44211, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 65, 5000, 5000, 5000, '$/docena de paquetes', 'Provincia de Cautín', 417, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44424, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 30, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44326, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 50, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44327, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 30, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44425, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 20, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44166, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 55, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44162, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 50, 5000, 5000, 5000, '$/docena de paquetes', 'Provincia de Cautín', 417, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44413, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 40, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44165, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 130, 5000, 6000, 5615, '$/docena de paquetes', 'Provincia de Cautín', 468, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44372, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 40, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44203, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 40, 5000, 5000, 5000, '$/docena de paquetes', 'Provincia de Cautín', 417, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44410, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 40, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44195, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 55, 5000, 5000, 5000, '$/docena de paquetes', 'Provincia de Cautín', 417, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44299, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 60, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44186, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 50, 5000, 5000, 5000, '$/docena de paquetes', 'Provincia de Cautín', 417, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44259, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 30, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44179, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 40, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44196, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 20, 5000, 5000, 5000, '$/docena de paquetes', 'Provincia de Cautín', 417, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44427, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 20, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44301, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 50, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44441, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 40, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44432, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 30, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44369, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 20, 4000, 4000, 4000, '$/docena de paquetes', 'Región Metropolitana', 333, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44302, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 20, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44411, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 20, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44329, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 40, 5000, 6000, 5500, '$/docena de paquetes', 'Provincia de Cautín', 458, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44438, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 30, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44428, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 10, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44442, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 20, 6000, 6000, 6000, '$/docena de paquetes', 'Provincia de Cautín', 500, 12, 'Hortaliza'),
        @(10, 'Vega Modelo de Temuco', 'La Araucanía', 44435, 9, 300000001, 'Rabanito', 'Sin especificar', 'Primera', 30, 7000, 7000, 7000, '$/docena de paquetes', 'Provincia de Cautín', 583, 12, 'Hortaliza')
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($c = 1; $c -le $vals.Count; $c++) {
        $ws.Cells.Item($row, $c).Value = $vals[$c - 1]
    }
}

# Rows 30-31 are brand new (the sheet only went to row 29 before); give the
# new "Fecha" cells (column D) the same date number format already used by
# the rest of the column so they match the existing table formatting.
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
$ws.Cells.Item(31, 4).NumberFormat = $ws.Cells.Item(2, 4).NumberFormat
